$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update K1 and A6 with new values (forced to text, matching existing
# numeric-looking string entries already used elsewhere in the sheet)
$ws.Range("K1").Value = "'9552"
$ws.Range("A6").Value = "'2824"
